# "edit link on ppt"
# Slide 11 / Shape 2 ("TextBox 3") holds the GitHub repository link.
# Update its position/size, enable explicit word-wrap, shrink the font,
# and swap in the new repository URL (while dropping the now-stale
# trailing endParaRPr that a plain Text= would otherwise keep around).

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(11)
$sh = $s.Shapes.Item(2)

# Turn on (already-default) word wrap explicitly so PowerPoint persists
# the attribute on <a:bodyPr wrap="square" .../>.
$sh.TextFrame.WordWrap = -1

# Reposition / resize the box (values are in points; PowerPoint COM
# stores coordinates as points and converts to EMU on save).
$sh.Left  = 99.12796
$sh.Width = 1232.8722

# Replace the URL. Delete()+InsertBefore() (rather than a plain
# TextRange.Text assignment) re-types the paragraph so the trailing
# <a:endParaRPr> produced by the old text goes away, while the run
# still inherits the existing color/typeface formatting.
$tr = $sh.TextFrame.TextRange
$tr.Delete()
$tr.InsertBefore("https://github.com/arawsardni/Final-Task---Home-Credit-Scorecard-Model/blob/main/Default_Prediction_Gaung_Taqwa.ipynb")
$tr.Font.Size = 24

# spAutoFit recalculates Height from the new text/width automatically;
# pin it to the final target height afterwards.
$sh.Height = 92.43394
